# Remove the no-longer-needed "D2" value ("STATE OF S TEXAS 3H1") from every
# worksheet in the workbook. Clearing the cell drops its shared-string
# reference, which in turn lets the now-unused string fall out of the
# shared strings table when the workbook is saved.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("D2").ClearContents()
}
